$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 629, shifting existing rows 629:670 down to 630:671.
$ws.Rows.Item(629).Insert()

# Force column A to be treated as text (not auto-converted to a date serial)
# before assigning the date-like string, matching the existing data format.
$ws.Range("A629").NumberFormat = "@"
$ws.Range("A629").Value = "2026/01/13"
$ws.Range("A629").Style = "Normal"

$ws.Range("B629").Value = "火"
$ws.Range("C629").Value = 20
$ws.Range("D629").Value = 201
